$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Logs" sheet: append the new test-mail row (row 48)
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 48

$logs.Cells.Item($newRow, 1).Value = "Wanneer komt mijn offerte?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #16: Wanneer komt mijn offerte?"
$logs.Cells.Item($newRow, 4).Value = "Offerte / Prijsaanvraag"
$logs.Cells.Item($newRow, 5).Value = "Geachte klant,`nDank u voor uw e-mail. Uw offerte zal naar verwachting binnen 24 uur worden verstuurd. Mocht u deze niet op tijd ontvangen, neem dan gerust contact met ons op.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item($newRow, 6).Value = "2025-06-26 23:53:56"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"

# Writing the multi-line "Antwoord" text makes the engine pin an explicit
# customHeight on the row (Excel parity for auto-measured rows); AutoFit
# clears that override again so row 48 stays on the sheet's default height,
# same as every other data row.
$logs.Rows.Item($newRow).AutoFit()

# Extend the conditional-formatting ranges (D/G/H/I) that covered 2:47 so
# they also cover the freshly added row 48. Modifying the AppliesTo range of
# one rule in a conditionalFormatting block re-targets the whole block.
$logs.Range("D2:D47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D48"))
$logs.Range("G2:G47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G48"))
$logs.Range("H2:H47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H48"))
$logs.Range("I2:I47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I48"))

# ---------------------------------------------------------------------------
# 2) "Dashboard" sheet: the new mail bumps the "Offerte / Prijsaanvraag"
#    category count from 2 to 3, re-sorting the summary table (rows 5-7).
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(5, 1).Value = "Offerte / Prijsaanvraag"
$dash.Cells.Item(5, 2).Value = 3

$dash.Cells.Item(6, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(6, 2).Value = 3

$dash.Cells.Item(7, 1).Value = "Overig"
$dash.Cells.Item(7, 2).Value = 3
